$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "3_3" (first sheet) - absolute (mA) measurements
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("3_3")

# Duplicate row 3's formatting down into row 4 (new "case_3" row) before we
# touch any values, so the new row inherits the same style / border as the
# existing data rows.
$ws1.Range("A3:G3").Copy($ws1.Range("A4:G4"))

# Row 2: case_0 -> Deep_Sleep, refreshed measurements
$ws1.Range("A2").Value = "Deep_Sleep"
$ws1.Range("B2").Value = -0.0011
$ws1.Range("C2").Value = 0.0092
$ws1.Range("D2").Value = -0.011
$ws1.Range("E2").Value = 0.005
$ws1.Range("F2").Value = 10
$ws1.Range("G2").Value = "-0.002798,-0.001482,-0.010965,0.009219,-0.001184,-0.000364,-0.000382,-0.005754,0.002791,-0.000463"

# Row 3: case_1, refreshed measurements
$ws1.Range("A3").Value = "case_1"
$ws1.Range("B3").Value = 0.0001
$ws1.Range("C3").Value = 0.0083
$ws1.Range("D3").Value = -0.0094
$ws1.Range("E3").Value = 0.0053
$ws1.Range("F3").Value = 10
$ws1.Range("G3").Value = "-0.004854,0.000177,-0.009424,-0.002978,-0.004835,0.005325,0.006064,0.002431,0.000673,0.008281"

# Row 4 (new): case_3
$ws1.Range("A4").Value = "case_3"
$ws1.Range("B4").Value = 0.0009
$ws1.Range("C4").Value = 0.0089
$ws1.Range("D4").Value = -0.0036
$ws1.Range("E4").Value = 0.004
$ws1.Range("F4").Value = 10
$ws1.Range("G4").Value = "-0.002699,-0.00178,0.008858,0.004188,0.000456,-0.00178,-0.002618,-0.003618,0.002133,0.006126"

# ---------------------------------------------------------------------------
# Sheet "1_8" (second sheet) - relative (%) measurements
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("1_8")

$ws2.Range("A3:G3").Copy($ws2.Range("A4:G4"))

# Row 2: case_0 -> Deep_Sleep, refreshed measurements
$ws2.Range("A2").Value = "Deep_Sleep"
$ws2.Range("B2").Value = -0.2405
$ws2.Range("C2").Value = -0.2352
$ws2.Range("D2").Value = -0.245
$ws2.Range("E2").Value = 0.0029
$ws2.Range("F2").Value = 10
$ws2.Range("G2").Value = "-0.240528,-0.244695,-0.240515,-0.244961,-0.236367,-0.235231,-0.240503,-0.241621,-0.239108,-0.241059"

# Row 3: case_1, refreshed measurements
$ws2.Range("A3").Value = "case_1"
$ws2.Range("B3").Value = -0.2416
$ws2.Range("C3").Value = -0.2377
$ws2.Range("D3").Value = -0.2497
$ws2.Range("E3").Value = 0.0038
$ws2.Range("F3").Value = 10
$ws2.Range("G3").Value = "-0.247745,-0.238861,-0.240509,-0.241349,-0.23799,-0.241102,-0.249684,-0.237737,-0.239682,-0.241355"

# Row 4 (new): case_3
$ws2.Range("A4").Value = "case_3"
$ws2.Range("B4").Value = -0.2423
$ws2.Range("C4").Value = -0.2366
$ws2.Range("D4").Value = -0.2483
$ws2.Range("E4").Value = 0.0037
$ws2.Range("F4").Value = 10
$ws2.Range("G4").Value = "-0.239694,-0.239991,-0.247473,-0.240225,-0.240799,-0.244992,-0.236644,-0.248295,-0.23996,-0.245251"
